# Burndown Chart updated for until today (12/02/17).
# - Fill in "Story Points (Done that day)" (row 2) for the days that have
#   now passed (columns F, G, H) with 0 points done.
# - Move the active selection to H3 to reflect the latest day reviewed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2:H2").Value = 0

$ws.Range("H3").Select()
